# Femacal de La Calera - Betarraga: add a new weekly price report.
#
# The sheet is a weekly time series (each reporting date contributes two
# rows: "Primera" and "Segunda" quality). A new, more recent date (serial
# 45265) is inserted at the top of the data block (just under the header
# rows), pushing every existing record down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 1213. This shifts
# the existing rows 1213:1304 down to 1215:1306 (and grows the used range
# from A1:R1304 to A1:R1306), while the new rows 1213:1214 inherit the
# date number-format from the row above (needed for column D).
$ws.Rows("1213:1214").Insert()

# ---- New row 1213 ("Primera" quality) ----
$ws.Range("A1213").Value = 3
$ws.Range("B1213").Value = "Femacal de La Calera"
$ws.Range("C1213").Value = "Coquimbo"
$ws.Range("D1213").Value = 45265
$ws.Range("E1213").Value = 5
$ws.Range("F1213").Value = 100114014
$ws.Range("G1213").Value = "Betarraga"
$ws.Range("H1213").Value = "Sin especificar"
$ws.Range("I1213").Value = "Primera"
$ws.Range("J1213").Value = 1800
$ws.Range("K1213").Value = 600
$ws.Range("L1213").Value = 600
$ws.Range("M1213").Value = 600
$ws.Range("N1213").Value = "`$/paquete 4 unidades"
$ws.Range("O1213").Value = "Provincia de Quillota"
$ws.Range("P1213").Value = 150
$ws.Range("Q1213").Value = 4
$ws.Range("R1213").Value = "Hortaliza"

# ---- New row 1214 ("Segunda" quality) ----
$ws.Range("A1214").Value = 3
$ws.Range("B1214").Value = "Femacal de La Calera"
$ws.Range("C1214").Value = "Coquimbo"
$ws.Range("D1214").Value = 45265
$ws.Range("E1214").Value = 5
$ws.Range("F1214").Value = 100114014
$ws.Range("G1214").Value = "Betarraga"
$ws.Range("H1214").Value = "Sin especificar"
$ws.Range("I1214").Value = "Segunda"
$ws.Range("J1214").Value = 1600
$ws.Range("K1214").Value = 400
$ws.Range("L1214").Value = 400
$ws.Range("M1214").Value = 400
$ws.Range("N1214").Value = "`$/paquete 4 unidades"
$ws.Range("O1214").Value = "Provincia de Quillota"
$ws.Range("P1214").Value = 100
$ws.Range("Q1214").Value = 4
$ws.Range("R1214").Value = "Hortaliza"
